$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric-looking text values in column D need a text-forcing trick so Excel
# keeps them as text (matching the source inlineStr cells) instead of coercing
# them to numbers.
function Set-TextValue($addr, $val) {
    $ws.Range($addr).Value = "'" + $val
    $ws.Range($addr).ClearFormats()
}

Set-TextValue 'D2' '247.77'
Set-TextValue 'E2' '1BNBBNB'
Set-TextValue 'D4' '5.554'
Set-TextValue 'D5' '0.05633'
Set-TextValue 'D8' '0.8021'
Set-TextValue 'D9' '1.068'
Set-TextValue 'D10' '0.1433'
Set-TextValue 'D11' '0.07398'
Set-TextValue 'D12' '0.03197'
Set-TextValue 'D13' '0.02971'
Set-TextValue 'D14' '0.09259'
Set-TextValue 'D15' '0.001675'
Set-TextValue 'B16' 'CoinExToken'
Set-TextValue 'C16' 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
Set-TextValue 'D16' '0.04723'
Set-TextValue 'E16' '15CoinExTokenCET'
Set-TextValue 'B17' 'One'
Set-TextValue 'C17' 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
Set-TextValue 'D17' '0.0005752'
Set-TextValue 'E17' '16OneONEWorstin24h'
Set-TextValue 'B18' 'TigerCash'
Set-TextValue 'C18' 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextValue 'D18' '0.006270'
Set-TextValue 'E18' '17TigerCashTCH'
Set-TextValue 'B19' 'BitKan'
Set-TextValue 'C19' 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
Set-TextValue 'D19' '0.001049'
Set-TextValue 'E19' '18BitKanKAN'
Set-TextValue 'B20' 'HotbitToken'
Set-TextValue 'C20' 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
Set-TextValue 'D20' '0.003823'
Set-TextValue 'E20' '19HotbitTokenHTB'
Set-TextValue 'B21' 'NitroEx'
Set-TextValue 'C21' 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
Set-TextValue 'D21' '0.0001501'
Set-TextValue 'E21' '20NitroExNTX'
Set-TextValue 'B22' 'UpBots'
Set-TextValue 'C22' 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
Set-TextValue 'D22' '0.0004602'
Set-TextValue 'E22' '21UpBotsUBXT'
Set-TextValue 'B23' 'LEO'
Set-TextValue 'C23' 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue 'D23' '3.981'
Set-TextValue 'E23' '22LEOLEO'
Set-TextValue 'B24' 'BTSEToken'
Set-TextValue 'C24' 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
Set-TextValue 'D24' '2.112'
Set-TextValue 'E24' '23BTSETokenBTSE'
Set-TextValue 'B25' 'BitpandaEcosystemToken'
Set-TextValue 'C25' 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
Set-TextValue 'D25' '0.3311'
Set-TextValue 'E25' '24BitpandaEcosystemTokenBEST'
Set-TextValue 'B26' 'ProBitToken'
Set-TextValue 'C26' 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
Set-TextValue 'D26' '0.1291'
Set-TextValue 'E26' '25ProBitTokenPROBBestin24h'
Set-TextValue 'B27' 'MCDex'
Set-TextValue 'C27' 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
Set-TextValue 'D27' '2.585'
Set-TextValue 'E27' '26MCDexMCB'
Set-TextValue 'D40' '0.04185'
Set-TextValue 'B41' 'KickToken'
Set-TextValue 'C41' 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
Set-TextValue 'D41' '0.007184'
Set-TextValue 'E41' '40KickTokenKICK'
Set-TextValue 'D42' '0.003502'
Set-TextValue 'B43' 'BKEXToken'
Set-TextValue 'C43' 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
Set-TextValue 'D43' '0.1045'
Set-TextValue 'E43' '42BKEXTokenBKK'
Set-TextValue 'D44' '0.008698'
Set-TextValue 'D45' '0.00005620'
Set-TextValue 'D48' '0.02785'
Set-TextValue 'D49' '0.00002101'
Set-TextValue 'D50' '0.01011'
